$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Neodymium")
$ws.Range("C2").Value = [double]"2.437496767148889E-05"
$ws.Range("D2").Value = 0.09088720653956146
$ws.Range("E2").Value = 0.4458999904664008
$ws.Range("B3").Value = [double]"2.183968937109617E-12"
$ws.Range("C3").Value = 0.001179565666721849
$ws.Range("D3").Value = 0.07925110730391915
$ws.Range("E3").Value = 0.3964840132176871
$ws.Range("B4").Value = [double]"3.409239115768776E-14"
$ws.Range("C4").Value = 0.001066828779339647
$ws.Range("D4").Value = 0.06477007120855503
$ws.Range("E4").Value = 0.3501040751224138
$ws.Range("C5").Value = [double]"2.374640474856044E-08"
$ws.Range("D5").Value = 0.003580656840928627
$ws.Range("E5").Value = 0.02923697393808194

$ws = $wb.Worksheets.Item("Dysprosium")
$ws.Range("C2").Value = [double]"2.76196144543578E-05"
$ws.Range("D2").Value = 0.07846518569028681
$ws.Range("E2").Value = 0.5052554730683479
$ws.Range("C3").Value = 0.001336582241976157
$ws.Range("D3").Value = 0.06841945184062941
$ws.Range("E3").Value = 0.4492615428244429
$ws.Range("C4").Value = 0.001208838508887113
$ws.Range("D4").Value = 0.05591761324890304
$ws.Range("E4").Value = 0.396707790718064
$ws.Range("C5").Value = [double]"2.690738107519781E-08"
$ws.Range("D5").Value = 0.003091270098551975
$ws.Range("E5").Value = 0.03312882129179109

$ws = $wb.Worksheets.Item("Copper")
$ws.Range("B2").Value = [double]"3.278472098474135E-06"
$ws.Range("C2").Value = 0.003032755898708039
$ws.Range("D2").Value = 0.5823644742193366
$ws.Range("E2").Value = 0.6659852986845426
$ws.Range("B3").Value = [double]"2.22924718813326E-05"
$ws.Range("C3").Value = 0.01094106145566737
$ws.Range("D3").Value = 0.4135478565326619
$ws.Range("E3").Value = 0.5113972578016854
$ws.Range("B4").Value = [double]"6.611256234481376E-05"
$ws.Range("C4").Value = 0.002928691932974728
$ws.Range("D4").Value = 0.3509623647012626
$ws.Range("E4").Value = 0.5150582967402831
$ws.Range("B5").Value = [double]"2.076903987060008E-05"
$ws.Range("C5").Value = 0.006421218662127675
$ws.Range("D5").Value = 0.5040422497805486
$ws.Range("E5").Value = 0.5201943099942047

$ws = $wb.Worksheets.Item("Raw silicon")
$ws.Range("B2").Value = [double]"4.966311329314386E-07"
$ws.Range("C2").Value = 0.0005182112305667888
$ws.Range("D2").Value = 0.3212156744753464
$ws.Range("E2").Value = 0.775899308654714
$ws.Range("B3").Value = [double]"5.299988190966853E-07"
$ws.Range("C3").Value = 0.001731461407493704
$ws.Range("D3").Value = 0.169168646809207
$ws.Range("E3").Value = 0.4311704658185009
$ws.Range("B4").Value = [double]"3.396088080967769E-06"
$ws.Range("C4").Value = 0.0004859890105311957
$ws.Range("D4").Value = 0.1748139439746624
$ws.Range("E4").Value = 0.5412910095236106
$ws.Range("B5").Value = [double]"1.82357813169359E-06"
$ws.Range("C5").Value = 0.0006171561129647565
$ws.Range("D5").Value = 0.2994238988798891
$ws.Range("E5").Value = 0.643701610727723
